# ---------------------------------------------------------------------------
# Adds a new "Player Info" worksheet (with the player's basic details) in
# front of the existing "ODI Batting" / "ODI Bowling" sheets, and reworks the
# MATCH_CARD_LINK column on both of those sheets into a MATCH_CODE column
# that stores just the numeric How Stat match code instead of the full URL.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Helper: write a text value into a cell without Excel coercing a
# numeric-looking string (e.g. "6245") into a real number. Purely-numeric
# text needs a "@" (Text) number format while it is assigned, otherwise
# Excel stores it as a real <v> number and drops leading zeros / turns it
# into a different type than the source workbook uses (inline/shared
# string). Re-applying the cell's own Style afterwards snaps the format
# back to the plain/un-decorated default, which is correct here because
# this helper is only ever used on plain (non-header, unstyled) data
# cells whose numeric-looking values need to stay text.
function Set-TextValue($cell, $text) {
    if ($text -match '^[0-9]+$') {
        $originalStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = $originalStyle
    } else {
        $cell.Value = $text
    }
}

# ---------------------------------------------------------------------------
# 1) Insert the new "Player Info" sheet before "ODI Batting".
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")

$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

# Worksheet handles captured before a sheet-collection structural change
# (Add/Delete/Move) are resolved positionally, so grab fresh handles by
# name now that "Player Info" occupies the old "ODI Batting" slot.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $headers.Length; $c++) {
    Set-TextValue $playerInfo.Cells.Item(1, $c) $headers[$c - 1]
}

$values = @("6245", "Christopher Nicholas Greaves", "Right Handed", "Right Arm Leg Break")
for ($c = 1; $c -le $values.Length; $c++) {
    Set-TextValue $playerInfo.Cells.Item(2, $c) $values[$c - 1]
}

# Match the bold / bordered / centered header look already used by the
# other sheets' header rows, by copying the real formatting over (copying
# values too would be fine, but PasteSpecial formats-only keeps our text).
$battingSheet.Range("A1:D1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)

$playerInfo.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2) ODI Batting: rename MATCH_CARD_LINK -> MATCH_CODE (column D) and turn
#    the full scorecard URL into just the trailing numeric match code.
# ---------------------------------------------------------------------------
Set-TextValue $battingSheet.Range("D1") "MATCH_CODE"

$lastRow = $battingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $url = $cell.Value()
    $code = $url.Substring($url.IndexOf("MatchCode=") + 10)
    Set-TextValue $cell $code
}

# ---------------------------------------------------------------------------
# 3) ODI Bowling: rename MATCH_CARD_LINK -> MATCH_CODE (column B) and turn
#    the full scorecard URL into just the trailing numeric match code.
# ---------------------------------------------------------------------------
Set-TextValue $bowlingSheet.Range("B1") "MATCH_CODE"

$lastRow = $bowlingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $url = $cell.Value()
    $code = $url.Substring($url.IndexOf("MatchCode=") + 10)
    Set-TextValue $cell $code
}
